# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" / "Latest Handoff Date" values for every
# file row that is not already fully synced / still in translation (rows 7, 10-16)
# with the timestamps produced by this handoff-generation run.

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

# zh-cn worksheet: "Latest Handoff Datetime" column E.
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhTimestamp = "2016-03-24 20:36:48"
foreach ($row in $rows) {
    $wsZh.Range("E$row").Value = $zhTimestamp
}

# de-de worksheet: "Latest Handoff Datetime" column E.
$wsDe = $wb.Worksheets.Item("de-de")
$deTimestamp = "2016-03-24 20:36:54"
foreach ($row in $rows) {
    $wsDe.Range("E$row").Value = $deTimestamp
}

# Overview worksheet: "Latest Handoff Date" mirrors the de-de handoff datetime.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in $rows) {
    $wsOverview.Range("D$row").Value = $deTimestamp
}
